$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values of rows 4, 5 and 6 across columns A:Q
$row4 = $ws.Range("A4:Q4").Value2
$row5 = $ws.Range("A5:Q5").Value2
$row6 = $ws.Range("A6:Q6").Value2

# Columns E (publication_date, e.g. "2022-07-21"), M (cited_by_count, e.g.
# "0") and N (publication_year, e.g. "2022") hold numeric/date-looking text
# that must stay plain text, not be auto-converted to a date serial / number.
# Pre-format them as text so the assignment below keeps them as strings.
$ws.Range("E4:E6").NumberFormat = "@"
$ws.Range("M4:M6").NumberFormat = "@"
$ws.Range("N4:N6").NumberFormat = "@"

# Rotate the rows: new row4 = old row6, new row5 = old row4, new row6 = old row5
$ws.Range("A4:Q4").Value2 = $row6
$ws.Range("A5:Q5").Value2 = $row4
$ws.Range("A6:Q6").Value2 = $row5

# Drop the temporary text formatting so the cells return to the default
# (unstyled) appearance, matching the rest of the sheet.
$ws.Range("E4:E6").ClearFormats()
$ws.Range("M4:M6").ClearFormats()
$ws.Range("N4:N6").ClearFormats()
